$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number format + style) from the last existing data column (BP) into the new column (BQ)
$ws.Range("BP2:BP11").Copy()
$ws.Range("BQ2:BQ11").PasteSpecial(-4122)

# Header for the new date column
$ws.Range("BQ1").Value = "14-sep"

# Data values for the new date column
$ws.Range("BQ2").Value = 14
$ws.Range("BQ3").Value = 13
$ws.Range("BQ4").Value = 12
$ws.Range("BQ5").Value = 12
$ws.Range("BQ6").Value = 10
$ws.Range("BQ7").Value = 17
$ws.Range("BQ8").Value = 20
$ws.Range("BQ9").Value = 9
$ws.Range("BQ10").Value = 7
$ws.Range("BQ11").Value = 11

# Match the saved selection state from the edited workbook
$ws.Range("BQ12").Select()
